# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matches the workbook's existing
# t="inlineStr"/shared-string cells) without leaving a lasting style/number-
# format change behind on the cell once we are done.
function Set-TextCell {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "29.427.58"
Set-TextCell $ws.Range("E2") "  -0.39%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.901.57"
Set-TextCell $ws.Range("E3") "  -0.62%  "

# Row 4
Set-TextCell $ws.Range("D4") "1.006"
Set-TextCell $ws.Range("E4") "  +0.57%  "

# Row 5
Set-TextCell $ws.Range("D5") "324.84"
Set-TextCell $ws.Range("E5") "  -1.22%  "

# Row 6
Set-TextCell $ws.Range("E6") "  +0.56%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.4828"
Set-TextCell $ws.Range("E7") "  +3.34%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.4063"
Set-TextCell $ws.Range("E8") "  -0.90%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.08139"
Set-TextCell $ws.Range("E9") "  +1.40%  "

# Row 10
Set-TextCell $ws.Range("E10") "  -0.73%  "

# Row 11
Set-TextCell $ws.Range("D11") "23.42"
Set-TextCell $ws.Range("E11") "  +4.31%  "

# Row 12
Set-TextCell $ws.Range("D12") "1.907.44"
Set-TextCell $ws.Range("E12") "  -0.39%  "

# Row 13
Set-TextCell $ws.Range("D13") "5.983"
Set-TextCell $ws.Range("E13") "  +0.35%  "

# Row 14
Set-TextCell $ws.Range("D14") "7.061"
Set-TextCell $ws.Range("E14") "  -1.33%  "

# Row 15
Set-TextCell $ws.Range("D15") "90.21"
Set-TextCell $ws.Range("E15") "  +0.73%  "

# Row 16
Set-TextCell $ws.Range("E16") "  +0.63%  "

# Row 17
Set-TextCell $ws.Range("D17") "0.06743"
Set-TextCell $ws.Range("E17") "  +2.47%  "

# Row 18
Set-TextCell $ws.Range("D18") "0.00001038"
Set-TextCell $ws.Range("E18") "  +0.78%  "

# Row 19
Set-TextCell $ws.Range("D19") "17.66"
Set-TextCell $ws.Range("E19") "  -0.25%  "

# Row 20
Set-TextCell $ws.Range("E20") "  +0.54%  "

# Row 21
Set-TextCell $ws.Range("D21") "29.443.40"
Set-TextCell $ws.Range("E21") "  -0.26%  "

# Row 22
Set-TextCell $ws.Range("D22") "5.556"
Set-TextCell $ws.Range("E22") "  +0.18%  "

# Row 23
Set-TextCell $ws.Range("D23") "11.81"
Set-TextCell $ws.Range("E23") "  +2.03%  "

# Row 24
Set-TextCell $ws.Range("D24") "2.158"
Set-TextCell $ws.Range("E24") "  -2.39%  "

# Row 25
Set-TextCell $ws.Range("D25") "2.140.85"
Set-TextCell $ws.Range("E25") "  -0.18%  "

# Row 26
Set-TextCell $ws.Range("D26") "154.09"
Set-TextCell $ws.Range("E26") "  +0.61%  "

# Row 27
Set-TextCell $ws.Range("D27") "20.00"
Set-TextCell $ws.Range("E27") "  +0.79%  "

# Row 28
Set-TextCell $ws.Range("D28") "6.152"
Set-TextCell $ws.Range("E28") "  +6.96%  "

# Row 29
Set-TextCell $ws.Range("D29") "2.089"
Set-TextCell $ws.Range("E29") "  -2.10%  "

# Row 30
Set-TextCell $ws.Range("D30") "118.90"
Set-TextCell $ws.Range("E30") "  +1.55%  "

# Row 31
Set-TextCell $ws.Range("D31") "1.029"
Set-TextCell $ws.Range("E31") "  -4.10%  "

# Row 32
Set-TextCell $ws.Range("D32") "0.09525"
Set-TextCell $ws.Range("E32") "  +0.65%  "

# Row 33
Set-TextCell $ws.Range("D33") "5.500"
Set-TextCell $ws.Range("E33") "  +1.96%  "

# Row 34
Set-TextCell $ws.Range("D34") "3.549"
Set-TextCell $ws.Range("E34") "  -0.64%  "

# Row 35
Set-TextCell $ws.Range("D35") "1.389"
Set-TextCell $ws.Range("E35") "  -2.80%  "

# Row 36
Set-TextCell $ws.Range("B36") "VeChain"
Set-TextCell $ws.Range("C36") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D36") "0.02257"
Set-TextCell $ws.Range("E36") "  -0.35%  "

# Row 37
Set-TextCell $ws.Range("B37") "Hedera"
Set-TextCell $ws.Range("C37") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D37") "0.06090"
Set-TextCell $ws.Range("E37") "  -0.21%  "

# Row 38
Set-TextCell $ws.Range("D38") "1.167"
Set-TextCell $ws.Range("E38") "  -0.96%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.5913"
Set-TextCell $ws.Range("E39") "  +0.44%  "

# Row 40
Set-TextCell $ws.Range("D40") "7.916"
Set-TextCell $ws.Range("E40") "  -6.26%  "

# Row 41
Set-TextCell $ws.Range("E41") "  +0.33%  "

# Row 42
Set-TextCell $ws.Range("E42") "  +0.64%  "

# Row 43
Set-TextCell $ws.Range("D43") "1.289"
Set-TextCell $ws.Range("E43") "  -1.44%  "

# Row 44
Set-TextCell $ws.Range("B44") "RenderToken"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D44") "2.410"
Set-TextCell $ws.Range("E44") "  +0.57%  "

# Row 45
Set-TextCell $ws.Range("B45") "Cronos"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws.Range("D45") "0.07719"
Set-TextCell $ws.Range("E45") "  +0.19%  "

# Row 46
Set-TextCell $ws.Range("D46") "12.31"
Set-TextCell $ws.Range("E46") "  +0.93%  "

# Row 47
Set-TextCell $ws.Range("D47") "0.5551"
Set-TextCell $ws.Range("E47") "  -0.25%  "

# Row 48
Set-TextCell $ws.Range("D48") "1.935"
Set-TextCell $ws.Range("E48") "  +0.47%  "

# Row 49
Set-TextCell $ws.Range("D49") "114.70"
Set-TextCell $ws.Range("E49") "  +1.08%  "

# Row 50
Set-TextCell $ws.Range("D50") "72.49"
Set-TextCell $ws.Range("E50") "  +1.56%  "

# Row 51
Set-TextCell $ws.Range("E51") "  +1.75%  "

